$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-5 (no longer present in the data)
$ws.Range("A3:B5").EntireRow.Delete()

# Update row 2 values: date serial and temperature
$ws.Range("A2").Value = 44835.125
$ws.Range("B2").Value = 12.4
